# Wind_Tunnel_Data.xlsx - add vertical tail (and fuselage / wing thickness / volume)
# characteristics to the "Airplane Characteristics" sheet, rename it to
# "Wing Characteristics" section header, and add supporting named ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Airplane Characteristics")

# ---------------------------------------------------------------------------
# 1. Make room: insert the new rows that land *inside* the pre-existing
#    block (rows below this just get written directly since they are past
#    the end of the original used range).
#    Row numbers below are the CURRENT (evolving) row index at the moment
#    of each insert - i.e. apply them top-to-bottom in this exact order.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Insert()   # new: Wing Thickness
$ws.Rows.Item(10).Insert()  # new: Wing Volume
$ws.Rows.Item(12).Insert()  # new: blank separator row
$ws.Rows.Item(13).Insert()  # new: "Horizontal Tail Characteristics" header
$ws.Rows.Item(17).Insert()  # new: Tail Thickness
$ws.Rows.Item(19).Insert()  # new: Tail Volume

# ---------------------------------------------------------------------------
# 2. Re-point the defined names whose target cell moved. Re-assigning
#    .RefersTo on the existing Name object leaves stale formula bindings in
#    this runtime, so delete + re-add instead.
# ---------------------------------------------------------------------------
$wb.Names.Item("Cchord").Delete()
$wb.Names.Item("CGw").Delete()
$wb.Names.Item("lamda_t").Delete()
$wb.Names.Item("lbar").Delete()
$wb.Names.Item("tamc").Delete()
$wb.Names.Item("tmac").Delete()
$wb.Names.Item("Trc").Delete()
$wb.Names.Item("Ttc").Delete()
$wb.Names.Item("wing_area").Delete()
$wb.Names.Item("wmac").Delete()

$wb.Names.Add("Cchord", "='Airplane Characteristics'!`$C`$8")
$wb.Names.Add("CGw", "='Airplane Characteristics'!`$C`$22")
$wb.Names.Add("lamda_t", "='Airplane Characteristics'!`$C`$16")
$wb.Names.Add("lbar", "='Airplane Characteristics'!`$C`$23")
$wb.Names.Add("tamc", "='Airplane Characteristics'!`$C`$21")
$wb.Names.Add("tmac", "='Airplane Characteristics'!`$C`$21")
$wb.Names.Add("Trc", "='Airplane Characteristics'!`$C`$14")
$wb.Names.Add("Ttc", "='Airplane Characteristics'!`$C`$15")
$wb.Names.Add("wing_area", "='Airplane Characteristics'!`$C`$9")
$wb.Names.Add("wmac", "='Airplane Characteristics'!`$C`$20")

# Brand-new defined names.
$wb.Names.Add("Bval", "='Airplane Characteristics'!`$G`$5")
$wb.Names.Add("flength", "='Airplane Characteristics'!`$C`$34")
$wb.Names.Add("h_vt", "='Airplane Characteristics'!`$C`$30")
$wb.Names.Add("l_vt", "='Airplane Characteristics'!`$C`$31")
$wb.Names.Add("lamda_vt", "='Airplane Characteristics'!`$C`$28")
$wb.Names.Add("t_tail", "='Airplane Characteristics'!`$C`$17")
$wb.Names.Add("t_wing", "='Airplane Characteristics'!`$C`$7")
$wb.Names.Add("tau1_", "='Airplane Characteristics'!`$G`$6")
$wb.Names.Add("vol_tail", "='Airplane Characteristics'!`$C`$19")
$wb.Names.Add("vol_wing", "='Airplane Characteristics'!`$C`$10")
$wb.Names.Add("vtMAC", "='Airplane Characteristics'!`$C`$29")
$wb.Names.Add("vtRc", "='Airplane Characteristics'!`$C`$26")
$wb.Names.Add("vtTc", "='Airplane Characteristics'!`$C`$27")

# ---------------------------------------------------------------------------
# 3. Column A header text & width (section renamed "Wing Characteristics").
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Wing Characteristics"
$ws.Columns.Item(1).ColumnWidth = 22.7109375

# ---------------------------------------------------------------------------
# 4. Wing block edits (rows 5-11).
# ---------------------------------------------------------------------------
$ws.Range("C5").Formula = "=24+3/8"

$ws.Range("A7").Value = "Wing Thickness (in)"
$ws.Range("B7").Value = "t_wing"
$ws.Range("C7").Value = 0.35
$ws.Range("C7").NumberFormat = "0.00"
$ws.Range("C7").Borders.Item(12).LineStyle = 1

$ws.Range("A10").Value = "Wing Volume (in^3)"
$ws.Range("B10").Value = "vol_wing"
$ws.Range("C10").Formula = "=wing_area*t_wing"
$ws.Range("C10").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 5. Horizontal Tail Characteristics section (rows 13-23).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Horizontal Tail Characteristics"

$ws.Range("A17").Value = "Tail Thickness"
$ws.Range("B17").Value = "t_tail"
$ws.Range("C17").Value = 0.17

$ws.Range("A19").Value = "Tail Volume (in^3)"
$ws.Range("B19").Value = "vol_tail"
$ws.Range("C19").Formula = "=C18*t_tail"
$ws.Range("C19").NumberFormat = "0.00"
$ws.Range("C19").Borders.Item(12).LineStyle = 1

# ---------------------------------------------------------------------------
# 6. Vertical Tail Characteristics section (rows 25-31, new).
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Vertical Tail Characteristics"

$ws.Range("A26").Value = "Root Chord (in)"
$ws.Range("B26").Value = "vtRc"
$ws.Range("C26").Value = 2.95

$ws.Range("A27").Value = "Tip Chord"
$ws.Range("B27").Value = "vtTc"
$ws.Range("C27").Value = 1.95

$ws.Range("A28").Value = "Taper Ratio"
$ws.Range("B28").Value = "lamda_vt"
$ws.Range("C28").Formula = "=vtTc/vtRc"
$ws.Range("C28").NumberFormat = "0.000"

$ws.Range("A29").Value = "MAC (in)"
$ws.Range("B29").Value = "vtMAC"
$ws.Range("C29").Formula = "=(2*vtRc)/3*((1+lamda_vt+lamda_vt^2)/(1+lamda_vt))"
$ws.Range("C29").NumberFormat = "0.000"

$ws.Range("A30").Value = "Height (in)"
$ws.Range("B30").Value = "h_vt"
$ws.Range("C30").Value = 2.15

$ws.Range("A31").Value = "Extrapolated Length (in)"
$ws.Range("B31").Value = "l_vt"
$ws.Range("C31").Value = 0.3

# ---------------------------------------------------------------------------
# 7. Fuselage Characteristics section (rows 33-37, new).
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "Fuselage Characteristics"
$ws.Range("A33:C33").Borders.Item(9).LineStyle = 1

$ws.Range("A34").Value = "Length (in)"
$ws.Range("B34").Value = "flength"
$ws.Range("C34").Value = 12

# Some trailing blank formatted rows/cells present in the final sheet.
$ws.Range("A35:C37").Value = ""
$ws.Range("A40:C40").Value = ""

# ---------------------------------------------------------------------------
# 8. Wind Tunnel block (columns E:G) - B Value & tau1 rows, and the new
#    "Correction Factors" sub-header underneath.
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "B Value "
$ws.Range("F5").Value = "Bval"
$ws.Range("G5").Formula = "=TunWidth/TunHeight"
$ws.Range("G5").NumberFormat = "0.000"

$ws.Range("E6").Value = "tau1"
$ws.Range("F6").Value = "tau1"
$ws.Range("G6").Value = 0.855

$ws.Range("E9").Value = "Correction Factors"
$ws.Range("E9:G9").Merge()

Write-Host "edit applied"
